# "changed footer not to be sticky, removed duplicate color"
#
# colors.xlsx lists swatches (hex code + R/G/B breakdown + light/dark
# foreground) one per row. Row 101 (E6E6E6) duplicated row 95, which was
# already E6E6E6 - remove the duplicate and let the rows below shift up.
#
# The sheet had also been left scrolled/selected far down (frozen-looking
# "sticky" footer area, K2:K106 selected, view scrolled to A91); reset the
# view to a normal, non-scrolled selection near the top of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the duplicate "E6E6E6" swatch row (row 101 duplicated row 95);
# everything below shifts up by one row and formulas/refs adjust automatically.
$ws.Rows("101").Delete() | Out-Null

# Move the selection/view back up so the sheet no longer opens scrolled to
# the bottom with a large sticky range selected.
$ws.Range("G10").Select() | Out-Null
